# Applies the wording/grammar corrections from the commit diff to
# '15-06-2017 - Cahier des charges.docx' using Find/Replace over the
# whole document (Word merges the matched + replacement text into a
# single run, same as a real Find & Replace All in the Word UI).
$d = $word.ActiveDocument

$d.Content.Find.Execute("’est dans cette objectif que ce projet comprendra", $true, $false, $false, $false, $false, $true, 1, $false, "’est dans cet objectif que ce projet comprendra", 2) | Out-Null
$d.Content.Find.Execute("leurs niveau d’acquisition des cours qu’ils ont reçût", $true, $false, $false, $false, $false, $true, 1, $false, "leurs niveaux d’acquisition des cours qu’ils ont reçus", 2) | Out-Null
$d.Content.Find.Execute("de leurs progression sur l’exercice en cours", $true, $false, $false, $false, $false, $true, 1, $false, "de leurs progressions sur l’exercice en cours", 2) | Out-Null
$d.Content.Find.Execute(" Toute ceci leurs permettra", $true, $false, $false, $false, $false, $true, 1, $false, " Tous ceci leurs permettra", 2) | Out-Null
$d.Content.Find.Execute("faciliter leurs suivi de recherches de stage", $true, $false, $false, $false, $false, $true, 1, $false, "faciliter leurs suivis de recherches de stage", 2) | Out-Null
$d.Content.Find.Execute("un système de couleur seras mis en place celons l’état de progression des échanges avec l’entreprise, les stagiaires pourrons aussi mettre les différent", $true, $false, $false, $false, $false, $true, 1, $false, "un système de couleur seras mis en place, selon l’état de progression des échanges avec l’entreprise, les stagiaires pourront aussi mettre les différent", 2) | Out-Null
$d.Content.Find.Execute("Une liste des entreprises qui recrute des stagiaires sera mise à disposition du formateur", $true, $false, $false, $false, $false, $true, 1, $false, "Une liste des entreprises qui recrutent des stagiaires sera mise à disposition du formateur", 2) | Out-Null
$d.Content.Find.Execute("qu’il y a des point à prendre en compte sur la plateforme", $true, $false, $false, $false, $false, $true, 1, $false, "qu’il y a des points à prendre en compte sur la plateforme", 2) | Out-Null
$d.Content.Find.Execute("Création de compte utilisateurs « ", $true, $false, $false, $false, $false, $true, 1, $false, "Création de compte utilisateur « ", 2) | Out-Null
$d.Content.Find.Execute("Système de notification ou d’envoi de mail automatique.", $true, $false, $false, $false, $false, $true, 1, $false, "Système de notification ou d’envois de mail automatique.", 2) | Out-Null
$d.Content.Find.Execute("ne sont pas attendues par le client a l’heure actuelle", $true, $false, $false, $false, $false, $true, 1, $false, "ne sont pas attendues par le client à l’heure actuelle", 2) | Out-Null
$d.Content.Find.Execute("elles restent toutefois lister dans ce documents pour de possible évolution", $true, $false, $false, $false, $false, $true, 1, $false, "elles restent toutefois lister dans ce document pour de possible évolution", 2) | Out-Null
$d.Content.Find.Execute("Saisie des feuilles pédagogique individuelle par le stagiaire", $true, $false, $false, $false, $false, $true, 1, $false, "Saisie des feuilles pédagogiques individuelle par le stagiaire", 2) | Out-Null
$d.Content.Find.Execute("jusqu’à validation des deux partie sans qu’aucune modification", $true, $false, $false, $false, $false, $true, 1, $false, "jusqu’à validation des deux parties sans qu’aucune modification", 2) | Out-Null
$d.Content.Find.Execute("trace des différentes modifications  apporter au document éditable", $true, $false, $false, $false, $false, $true, 1, $false, "trace des différentes modifications  apportée au document éditable", 2) | Out-Null
$d.Content.Find.Execute("tableau de niveau d’acquisition de l’apprentissage, et sur son ressenti", $true, $false, $false, $false, $false, $true, 1, $false, "tableau de niveau d’acquisition de l’apprentissage et sur son ressenti", 2) | Out-Null
$d.Content.Find.Execute("sur son avancement personnelle", $true, $false, $false, $false, $false, $true, 1, $false, "sur son avancement personnel", 2) | Out-Null
$d.Content.Find.Execute("Saisie des comptes rendu des réunions pédagogique avec édition", $true, $false, $false, $false, $false, $true, 1, $false, "Saisie des comptes rendus des réunions pédagogiques avec édition", 2) | Out-Null
$d.Content.Find.Execute("toutes les informations qui sera lié au compte stagiaire", $true, $false, $false, $false, $false, $true, 1, $false, "toutes les informations qui seront lié au compte stagiaire", 2) | Out-Null
$d.Content.Find.Execute("Suivi libre sur tous les thèmes lié à la formation (recherche de stage, difficulté personne etc…)", $true, $false, $false, $false, $false, $true, 1, $false, "Suivi libre sur tous les thèmes liés à la formation (recherche de stage, difficulté personnelle etc…)", 2) | Out-Null
$d.Content.Find.Execute("Trois types de compte différent « Formateur », « Stagiaire », «Collaborateur», qui aurons chacun des droits différent (non définit", $true, $false, $false, $false, $false, $true, 1, $false, "Trois types de compte différent « Formateur », « Stagiaire », « Collaborateur », qui auront chacun des droits différents (non définit", 2) | Out-Null
$d.Content.Find.Execute("le login pour les comptes stagiaire sera leur code OSIA", $true, $false, $false, $false, $false, $true, 1, $false, "le login pour les comptes stagiaires sera leur code OSIA", 2) | Out-Null
$d.Content.Find.Execute("liste des documents pas encore définît a ce jour", $true, $false, $false, $false, $false, $true, 1, $false, "liste des documents pas encore définît à ce jour", 2) | Out-Null
$d.Content.Find.Execute("Formulaire d’absence stagiaire avec sélections de la raison celons les codes déjà en vigueur", $true, $false, $false, $false, $false, $true, 1, $false, "Formulaire d’absence stagiaire avec sélections de la raison, selon les codes déjà en vigueur", 2) | Out-Null
$d.Content.Find.Execute("(catégorie à définir celons document Excel suivit de stage existant)", $true, $false, $false, $false, $false, $true, 1, $false, "(catégorie à définir, selon document Excel suivi de stage existant)", 2) | Out-Null
$d.Content.Find.Execute("Possibilité de publier son C.V, lettre de motivation utilisé lors des différentes recherches de stage", $true, $false, $false, $false, $false, $true, 1, $false, "Possibilité de publier son C.V, lettre de motivation utilisée lors des différentes recherches de stage", 2) | Out-Null
$d.Content.Find.Execute("Liste des entreprises ayant déjà pris des stagiaires, tenu à jour par le formateur, qui pourras être mis", $true, $false, $false, $false, $false, $true, 1, $false, "Liste des entreprises ayant déjà pris des stagiaires, tenu à jour par le formateur, qui pourra être mis", 2) | Out-Null
$d.Content.Find.Execute("ce système de notification ce feras ou directement sur le site", $true, $false, $false, $false, $false, $true, 1, $false, "ce système de notification se feras ou directement sur le site", 2) | Out-Null
$d.Content.Find.Execute("La couleur principal du site sera le bleu technologique, pour rappeler le monde informatique.", $true, $false, $false, $false, $false, $true, 1, $false, "La couleur principale du site sera le bleu.", 2) | Out-Null
$d.Content.Find.Execute("trop administratif a était conseiller, il a aussi était demander de rendre l’applicatif le plus ergonomique possible, et le plus simple possible.", $true, $false, $false, $false, $false, $true, 1, $false, "trop administratif a était conseillé, il a aussi était demander de rendre l’applicatif le plus ergonomique possible et le plus simple possible.", 2) | Out-Null
$d.Content.Find.Execute("aucune demande particulière n’a était formulé, donc le style appliqué au stagiaire sera aussi utiliser pour les partie formateur.", $true, $false, $false, $false, $false, $true, 1, $false, "aucune demande particulière n’a été formulé, donc le style appliqué au stagiaire sera aussi utilisé pour les parties formateur.", 2) | Out-Null

# Word leaves a zero-length "_GoBack" bookmark at the last edit location;
# put it where the diff shows it, right after the word that became "utilise".
$goBack = $d.Content
$goBack.Find.Execute("sera aussi utilisé", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBack) | Out-Null

